$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row
# (rows 2-412). The commit bumps that date by one day (46081 -> 46082)
# for every row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 412 }

$range = $ws.Range("C2:C$lastRow")
foreach ($cell in $range.Cells) {
    $cell.Value2 = $cell.Value2 + 1
}
